$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update the "Förändrad" (Changed) date column for rows 2-8
# from serial 46072 (2026-02-19) to serial 46073 (2026-02-20)
for ($row = 2; $row -le 8; $row++) {
    $ws.Cells.Item($row, 3).Value = 46073
}
